# Scheduled-runner update: refresh cached market-board pricing / profit
# figures on the per-job "Hades_Profits" sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Only literal value cells are touched -- no formulas,
# formatting, or structure changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 93.333336
$ws.Range("I2").Value = 92.5
$ws.Range("K2").Value = 92.5
$ws.Range("M2").Value = 20.5
# Row 100
$ws.Range("H100").Value = 2903.6365
$ws.Range("I100").Value = 2668
$ws.Range("J100").Value = 3100
$ws.Range("K100").Value = 2668
$ws.Range("L100").Value = 3100
$ws.Range("M100").Value = -2127
$ws.Range("N100").Value = -4182
# Row 129
$ws.Range("H129").Value = 926.8182
$ws.Range("I129").Value = 552.8333
$ws.Range("J129").Value = 1067.0625
$ws.Range("K129").Value = 1658.4999
$ws.Range("L129").Value = 3201.1875
$ws.Range("M129").Value = 3341.5001
$ws.Range("N129").Value = -13201.1875
# Row 137
$ws.Range("H137").Value = 2003649
$ws.Range("I137").Value = 3031380.5
$ws.Range("J137").Value = 8640.706
$ws.Range("K137").Value = 9094141.5
$ws.Range("L137").Value = 25922.118
$ws.Range("M137").Value = -9091591.5
$ws.Range("N137").Value = -31022.118
# Row 138
$ws.Range("H138").Value = 5715960
$ws.Range("I138").Value = 1305.7391
$ws.Range("J138").Value = 16669048
$ws.Range("K138").Value = 3917.2173
$ws.Range("L138").Value = 50007144
$ws.Range("M138").Value = 1222.7827
$ws.Range("N138").Value = -50017424
# Row 139
$ws.Range("H139").Value = 56500
$ws.Range("J139").Value = 56500
$ws.Range("L139").Value = 56500
$ws.Range("N139").Value = -66780

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1739.6
$ws.Range("I45").Value = 1649.5
$ws.Range("J45").Value = 1799.6666
$ws.Range("K45").Value = 1649.5
$ws.Range("L45").Value = 1799.6666
$ws.Range("M45").Value = -1272.5
$ws.Range("N45").Value = -2553.6666
# Row 97
$ws.Range("H97").Value = 4465287
$ws.Range("I97").Value = 5682931
$ws.Range("J97").Value = 593
$ws.Range("K97").Value = 5682931
$ws.Range("L97").Value = 593
$ws.Range("M97").Value = -5682435
$ws.Range("N97").Value = -1585
# Row 110
$ws.Range("H110").Value = 1065.8462
$ws.Range("I110").Value = 918.2727
$ws.Range("J110").Value = 1877.5
$ws.Range("K110").Value = 918.2727
$ws.Range("L110").Value = 1877.5
$ws.Range("M110").Value = 1126.7273
$ws.Range("N110").Value = -5967.5
# Row 132
$ws.Range("H132").Value = 47440.24
$ws.Range("I132").Value = 33010.594
$ws.Range("J132").Value = 80422.28999999999
$ws.Range("K132").Value = 99031.78199999999
$ws.Range("L132").Value = 241266.87
$ws.Range("M132").Value = -96501.78199999999
$ws.Range("N132").Value = -246326.87

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1125.0625
$ws.Range("I99").Value = 1017.2727
$ws.Range("J99").Value = 1362.2
$ws.Range("K99").Value = 1017.2727
$ws.Range("L99").Value = 1362.2
$ws.Range("M99").Value = 480.7273
$ws.Range("N99").Value = -4358.2
# Row 107
$ws.Range("H107").Value = 2123.25
$ws.Range("I107").Value = 1937.8636
$ws.Range("J107").Value = 2531.1
$ws.Range("K107").Value = 1937.8636
$ws.Range("L107").Value = 2531.1
$ws.Range("M107").Value = -17.86359999999991
$ws.Range("N107").Value = -6371.1

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2101.3635
$ws.Range("I16").Value = 2086
$ws.Range("J16").Value = 2128.25
$ws.Range("K16").Value = 2086
$ws.Range("L16").Value = 2128.25
$ws.Range("M16").Value = -1799
$ws.Range("N16").Value = -2702.25
# Row 94
$ws.Range("H94").Value = 4383.8887
$ws.Range("I94").Value = 9160
$ws.Range("J94").Value = 1344.5454
$ws.Range("K94").Value = 9160
$ws.Range("L94").Value = 1344.5454
$ws.Range("M94").Value = -8709
$ws.Range("N94").Value = -2246.5454
# Row 113
$ws.Range("H113").Value = 2101.3635
$ws.Range("I113").Value = 2086
$ws.Range("J113").Value = 2128.25
$ws.Range("K113").Value = 2086
$ws.Range("L113").Value = 2128.25
$ws.Range("M113").Value = 84
$ws.Range("N113").Value = -6468.25
# Row 121
$ws.Range("H121").Value = 30000
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
# Row 132
$ws.Range("H132").Value = 46913.523
$ws.Range("I132").Value = 3334.8572
$ws.Range("J132").Value = 114702.555
$ws.Range("K132").Value = 10004.5716
$ws.Range("L132").Value = 344107.665
$ws.Range("M132").Value = -7474.571599999999
$ws.Range("N132").Value = -349167.665

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 564.8039
$ws.Range("I107").Value = 486.7586
$ws.Range("J107").Value = 667.6818
$ws.Range("K107").Value = 1460.2758
$ws.Range("L107").Value = 2003.0454
$ws.Range("M107").Value = 459.7242000000001
$ws.Range("N107").Value = -5843.0454
# Row 118
$ws.Range("H118").Value = 4504.4546
$ws.Range("I118").Value = 825
$ws.Range("J118").Value = 5322.1113
$ws.Range("K118").Value = 2475
$ws.Range("L118").Value = 15966.3339
$ws.Range("M118").Value = -1232
$ws.Range("N118").Value = -18452.3339
# Row 132
$ws.Range("H132").Value = 3472.25
$ws.Range("I132").Value = 905.8
$ws.Range("J132").Value = 5305.4287
$ws.Range("K132").Value = 8152.2
$ws.Range("L132").Value = 47748.85830000001
$ws.Range("M132").Value = -5622.2
$ws.Range("N132").Value = -52808.85830000001
# Row 140
$ws.Range("H140").Value = 2204.3015
$ws.Range("I140").Value = 2862.85
$ws.Range("K140").Value = 8588.549999999999
$ws.Range("M140").Value = -3408.549999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 2365
$ws.Range("I23").Value = 230
$ws.Range("K23").Value = 230
$ws.Range("M23").Value = -7
# Row 97
$ws.Range("H97").Value = 2240
$ws.Range("I97").Value = 2308.889
$ws.Range("J97").Value = 2033.3334
$ws.Range("K97").Value = 2308.889
$ws.Range("L97").Value = 2033.3334
$ws.Range("M97").Value = -1812.889
$ws.Range("N97").Value = -3025.3334
# Row 107
$ws.Range("H107").Value = 599.5238000000001
$ws.Range("I107").Value = 224.27777
$ws.Range("J107").Value = 2851
$ws.Range("K107").Value = 224.27777
$ws.Range("L107").Value = 2851
$ws.Range("M107").Value = 1695.72223
$ws.Range("N107").Value = -6691
# Row 122
$ws.Range("H122").Value = 1733.6
$ws.Range("I122").Value = 1602
$ws.Range("J122").Value = 2040.6666
$ws.Range("K122").Value = 4806
$ws.Range("L122").Value = 6121.9998
$ws.Range("M122").Value = -2356
$ws.Range("N122").Value = -11021.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1688.5
$ws.Range("I68").Value = 1672.1052
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1672.1052
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -923.1052
$ws.Range("N68").Value = -3498
# Row 71
$ws.Range("H71").Value = 1688.5
$ws.Range("I71").Value = 1672.1052
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 8360.526
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -4616.526
$ws.Range("N71").Value = -17488
# Row 132
$ws.Range("H132").Value = 58549.277
$ws.Range("I132").Value = 2115.9167
$ws.Range("J132").Value = 171416
$ws.Range("K132").Value = 6347.750100000001
$ws.Range("L132").Value = 514248
$ws.Range("M132").Value = -3817.750100000001
$ws.Range("N132").Value = -519308

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 102
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
